# "Se agrego My Order View"
# Adds two new worksheets ("my_order" and "near_me") with device-measure
# tables (mirroring the existing restaurant / retaurant_proudct sheets),
# converts the restaurant sheet's per-row % formulas into one shared
# formula, and updates the active-tab/selection state accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) restaurant sheet: C4:C9 "B/B$1*100" formulas -> single shared formula
# ---------------------------------------------------------------------
$restaurant = $wb.Worksheets.Item("restaurant")
$restaurant.Range("C4:C9").Formula = "=B4/B`$1*100"

# ---------------------------------------------------------------------
# 2) New sheet "my_order" (sheetId 4), inserted right after
#    "retaurant_proudct"
# ---------------------------------------------------------------------
$afterMyOrder = $wb.Worksheets.Item($wb.Worksheets.Count)
$myOrder = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterMyOrder)
$myOrder.Name = "my_order"

$myOrder.Range("B1").Value = "alto"
$myOrder.Range("C1").Value = "ancho"

$myOrder.Range("A2").Value = "dispositivo"
$myOrder.Range("B2").Value = 385
$myOrder.Range("C2").Value = 178

$myOrder.Range("A3").Value = "buttons"
$myOrder.Range("B3").Value = 19
$myOrder.Range("C3").Value = 70

$myOrder.Range("B4").Formula = "=B3/B2*100"
$myOrder.Range("C4").Formula = "=C3/C2*100"

$myOrder.Range("A2:C2").Font.Bold = $true
$myOrder.Range("A3:C4").Interior.Color = 65535

[void]$myOrder.Range("A1:C4").Select()
$myOrder.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 3) New sheet "near_me" (sheetId 5), inserted right after "my_order" -
#    becomes the active / selected tab
# ---------------------------------------------------------------------
$afterNearMe = $wb.Worksheets.Item($wb.Worksheets.Count)
$nearMe = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterNearMe)
$nearMe.Name = "near_me"

$nearMe.Range("B1").Value = "alto"
$nearMe.Range("C1").Value = "ancho"

$nearMe.Range("A2").Value = "dispositivo"
$nearMe.Range("B2").Value = 460
$nearMe.Range("C2").Value = 212

$nearMe.Range("A3").Value = "title"
$nearMe.Range("B3").Value = 17
$nearMe.Range("C3").Value = 128

$nearMe.Range("B4").Formula = "=B3/B2*100"
$nearMe.Range("C4").Formula = "=C3/C2*100"

$nearMe.Range("A2:C2").Font.Bold = $true
$nearMe.Range("A3:C4").Interior.Color = 65535

[void]$nearMe.Range("B4").Select()

Write-Host "my_order / near_me sheets added"
